$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row for the new "peptide range" table (row 20)
$ws.Range("G20").Value = "sample"
$ws.Range("H20").Value = "peptide"
$ws.Range("I20").Value = "start_seq"
$ws.Range("J20").Value = "end_seq"
$ws.Range("K20").Value = "peptide_range"

# Copy the peptide name / sequence table (G6:H17) down to G21:H32
$srcRow = 6
for ($r = 21; $r -le 32; $r++) {
    $ws.Range("G$r").Formula = $ws.Range("G$srcRow").Formula
    $ws.Range("H$r").Formula = $ws.Range("H$srcRow").Formula
    $srcRow++
}

# Column I: start_seq, always 1
$ws.Range("I21:I32").Formula = "1"

# Column J: end_seq = LEN(H) - first cell own formula, rest as a shared-formula block
$ws.Range("J21").Formula = "=LEN(H21)"
$ws.Range("J22:J32").Formula = "=LEN(H22)"

# Column K: peptide_range = CONCAT(TEXT(I,"0000"),"-",TEXT(J,"0000"))
$ws.Range("K21").Formula = "=CONCAT(TEXT(I21,""0000""),""-"",TEXT(J21,""0000""))"
$ws.Range("K22:K32").Formula = "=CONCAT(TEXT(I22,""0000""),""-"",TEXT(J22,""0000""))"

# Scroll / selection state to match the saved view
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("L22").Select()
